$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 840.2222
$ws.Range("I43").Value = 1100
$ws.Range("J43").Value = 807.75
$ws.Range("K43").Value = 1100
$ws.Range("L43").Value = 807.75
$ws.Range("M43").Value = -1031
$ws.Range("N43").Value = -945.75

$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2402.2
$ws.Range("I2").Value = 2037
$ws.Range("K2").Value = 2037
$ws.Range("M2").Value = -1924

$ws.Range("H45").Value = 1451.3846
$ws.Range("I45").Value = 1439.0834
$ws.Range("K45").Value = 1439.0834
$ws.Range("M45").Value = -1062.0834

$ws.Range("H110").Value = 2367.889
$ws.Range("I110").Value = 2473
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 2473
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = -428
$ws.Range("N110").Value = -6090

$ws.Range("H116").Value = 2402.2
$ws.Range("I116").Value = 2037
$ws.Range("K116").Value = 2037
$ws.Range("M116").Value = 257

$ws.Range("H132").Value = 1928.0541
$ws.Range("I132").Value = 1357.5333
$ws.Range("K132").Value = 4072.5999
$ws.Range("M132").Value = -1542.5999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2402.2
$ws.Range("I3").Value = 2037
$ws.Range("K3").Value = 2037
$ws.Range("M3").Value = -1923

$ws.Range("H107").Value = 2589.361
$ws.Range("I107").Value = 1792.0344
$ws.Range("J107").Value = 5892.5713
$ws.Range("K107").Value = 1792.0344
$ws.Range("L107").Value = 5892.5713
$ws.Range("M107").Value = 127.9656
$ws.Range("N107").Value = -9732.5713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1023.9
$ws.Range("I16").Value = 1009.2
$ws.Range("K16").Value = 1009.2
$ws.Range("M16").Value = -722.2

$ws.Range("H31").Value = 1592.0312
$ws.Range("J31").Value = 2248.8
$ws.Range("L31").Value = 2248.8
$ws.Range("N31").Value = -2838.8

$ws.Range("H34").Value = 1592.0312
$ws.Range("J34").Value = 2248.8
$ws.Range("L34").Value = 2248.8
$ws.Range("N34").Value = -2652.8

$ws.Range("H113").Value = 1023.9
$ws.Range("I113").Value = 1009.2
$ws.Range("K113").Value = 1009.2
$ws.Range("M113").Value = 1160.8

$ws.Range("H132").Value = 3263.9697
$ws.Range("I132").Value = 2944.4075
$ws.Range("K132").Value = 8833.2225
$ws.Range("M132").Value = -6303.2225

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 386.42426
$ws.Range("I5").Value = 302.3846
$ws.Range("J5").Value = 698.5714
$ws.Range("K5").Value = 907.1537999999999
$ws.Range("L5").Value = 2095.7142
$ws.Range("M5").Value = -795.1537999999999
$ws.Range("N5").Value = -2319.7142

$ws.Range("H135").Value = 386.42426
$ws.Range("I135").Value = 302.3846
$ws.Range("J135").Value = 698.5714
$ws.Range("K135").Value = 2721.4614
$ws.Range("L135").Value = 6287.1426
$ws.Range("M135").Value = -186.4613999999997
$ws.Range("N135").Value = -11357.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 92.666664
$ws.Range("I2").Value = 81.2
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 81.2
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = 31.8
$ws.Range("N2").Value = -376

$ws.Range("H18").Value = 5000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 5000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5586
$ws.Range("M18").ClearContents()

$ws.Range("H43").Value = 5538.6665
$ws.Range("I43").Value = 1016
$ws.Range("J43").Value = 7800
$ws.Range("K43").Value = 1016
$ws.Range("L43").Value = 7800
$ws.Range("M43").Value = -865
$ws.Range("N43").Value = -8102

$ws.Range("H46").Value = 23333.334
$ws.Range("J46").Value = 23333.334
$ws.Range("L46").Value = 23333.334
$ws.Range("N46").Value = -23645.334

$ws.Range("H57").Value = 13295.818
$ws.Range("I57").Value = 6351.6665
$ws.Range("J57").Value = 15899.875
$ws.Range("K57").Value = 6351.6665
$ws.Range("L57").Value = 15899.875
$ws.Range("M57").Value = -5531.6665
$ws.Range("N57").Value = -17539.875

$ws.Range("H80").Value = 2827.1924
$ws.Range("I80").Value = 2714.75
$ws.Range("J80").Value = 3202
$ws.Range("K80").Value = 2714.75
$ws.Range("L80").Value = 3202
$ws.Range("M80").Value = -1716.75
$ws.Range("N80").Value = -5198

$ws.Range("H83").Value = 2827.1924
$ws.Range("I83").Value = 2714.75
$ws.Range("J83").Value = 3202
$ws.Range("K83").Value = 13573.75
$ws.Range("L83").Value = 16010
$ws.Range("M83").Value = -8581.75
$ws.Range("N83").Value = -25994

$ws.Range("H113").Value = 1622.0968
$ws.Range("J113").Value = 2164.3333
$ws.Range("L113").Value = 2164.3333
$ws.Range("N113").Value = -6504.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1520
$ws.Range("I81").Value = 1533.3334
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 3066.6668
$ws.Range("L81").Value = 3000
$ws.Range("M81").Value = -2005.6668
$ws.Range("N81").Value = -5122

$ws.Range("H84").Value = 1520
$ws.Range("I84").Value = 1533.3334
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 15333.334
$ws.Range("L84").Value = 15000
$ws.Range("M84").Value = -10029.334
$ws.Range("N84").Value = -25608

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H126").Value = 1809.4445
$ws.Range("I126").Value = 1031.25
$ws.Range("J126").Value = 2432
$ws.Range("K126").Value = 3093.75
$ws.Range("L126").Value = 7296
$ws.Range("M126").Value = -623.75
$ws.Range("N126").Value = -12236
